$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 500950
$ws.Range("I12").Value = 900
$ws.Range("J12").Value = 1001000
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 1001000
$ws.Range("M12").Value = -730
$ws.Range("N12").Value = -1001340

$ws.Range("H32").Value = 475
$ws.Range("I32").Value = 300
$ws.Range("J32").Value = 533.3333
$ws.Range("K32").Value = 300
$ws.Range("L32").Value = 533.3333
$ws.Range("M32").Value = 26
$ws.Range("N32").Value = -1185.3333

$ws.Range("H38").Value = 1101.7084
$ws.Range("I38").Value = 153.8125
$ws.Range("J38").Value = 2997.5
$ws.Range("K38").Value = 461.4375
$ws.Range("L38").Value = 8992.5
$ws.Range("M38").Value = -89.4375
$ws.Range("N38").Value = -9736.5

$ws.Range("H87").Value = 28420.88
$ws.Range("J87").Value = 28420.88
$ws.Range("L87").Value = 28420.88
$ws.Range("N87").Value = -30916.88

$ws.Range("H90").Value = 28420.88
$ws.Range("J90").Value = 28420.88
$ws.Range("L90").Value = 85262.64
$ws.Range("N90").Value = -97742.64

$ws.Range("H97").Value = 6195
$ws.Range("J97").Value = 6195
$ws.Range("L97").Value = 18585
$ws.Range("N97").Value = -19577

$ws.Range("H98").Value = 3855.4375
$ws.Range("I98").Value = 4322.077
$ws.Range("J98").Value = 1833.3334
$ws.Range("K98").Value = 4322.077
$ws.Range("L98").Value = 1833.3334
$ws.Range("M98").Value = -2824.077
$ws.Range("N98").Value = -4829.3334

$ws.Range("H103").Value = 425.4737
$ws.Range("I103").Value = 305
$ws.Range("J103").Value = 457.6
$ws.Range("K103").Value = 915
$ws.Range("L103").Value = 1372.8
$ws.Range("M103").Value = -329
$ws.Range("N103").Value = -2544.8

$ws.Range("H106").Value = 2950
$ws.Range("I106").Value = 500
$ws.Range("J106").Value = 3766.6667
$ws.Range("K106").Value = 500
$ws.Range("L106").Value = 3766.6667
$ws.Range("M106").Value = 131
$ws.Range("N106").Value = -5028.6667

$ws.Range("H122").Value = 3855.4375
$ws.Range("I122").Value = 4322.077
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 12966.231
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -10516.231
$ws.Range("N122").Value = -10400.0002

$ws.Range("H138").Value = 1928.5333
$ws.Range("I138").Value = 1155.0476
$ws.Range("J138").Value = 3733.3333
$ws.Range("K138").Value = 3465.142800000001
$ws.Range("L138").Value = 11199.9999
$ws.Range("M138").Value = 1674.857199999999
$ws.Range("N138").Value = -21479.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 29222.223
$ws.Range("J64").Value = 29222.223
$ws.Range("L64").Value = 29222.223
$ws.Range("N64").Value = -29718.223

$ws.Range("H67").Value = 29222.223
$ws.Range("J67").Value = 29222.223
$ws.Range("L67").Value = 29222.223
$ws.Range("N67").Value = -30938.223

$ws.Range("H74").Value = 1607.4615
$ws.Range("I74").Value = 1141.9375
$ws.Range("J74").Value = 2352.3
$ws.Range("K74").Value = 1141.9375
$ws.Range("L74").Value = 2352.3
$ws.Range("M74").Value = -267.9375
$ws.Range("N74").Value = -4100.3

$ws.Range("H77").Value = 1607.4615
$ws.Range("I77").Value = 1141.9375
$ws.Range("J77").Value = 2352.3
$ws.Range("K77").Value = 5709.6875
$ws.Range("L77").Value = 11761.5
$ws.Range("M77").Value = -1341.6875
$ws.Range("N77").Value = -20497.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4026
$ws.Range("I31").Value = 2852.6
$ws.Range("J31").Value = 5981.6665
$ws.Range("K31").Value = 2852.6
$ws.Range("L31").Value = 5981.6665
$ws.Range("M31").Value = -2557.6
$ws.Range("N31").Value = -6571.6665

$ws.Range("H34").Value = 4026
$ws.Range("I34").Value = 2852.6
$ws.Range("J34").Value = 5981.6665
$ws.Range("K34").Value = 2852.6
$ws.Range("L34").Value = 5981.6665
$ws.Range("M34").Value = -2650.6
$ws.Range("N34").Value = -6385.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 276.02777
$ws.Range("I5").Value = 276.02777
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 828.08331
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -716.08331
$ws.Range("N5").ClearContents()

$ws.Range("H52").Value = 731.3333
$ws.Range("J52").Value = 731.3333
$ws.Range("L52").Value = 2193.9999
$ws.Range("N52").Value = -2725.9999

$ws.Range("H59").Value = 2501.2
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2501.2
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 7503.599999999999
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -8583.599999999999

$ws.Range("H81").Value = 14978.75
$ws.Range("I81").Value = 277
$ws.Range("K81").Value = 831
$ws.Range("M81").Value = 292

$ws.Range("H84").Value = 14978.75
$ws.Range("I84").Value = 277
$ws.Range("K84").Value = 2493
$ws.Range("M84").Value = 3123

$ws.Range("H104").Value = 2961.48
$ws.Range("I104").Value = 2700
$ws.Range("J104").Value = 2972.375
$ws.Range("K104").Value = 8100
$ws.Range("L104").Value = 8917.125
$ws.Range("M104").Value = -5479
$ws.Range("N104").Value = -14159.125

$ws.Range("H105").Value = 4900
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 743.36365
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 743.36365
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2230.09095
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6570.09095

$ws.Range("H118").Value = 1704.1428
$ws.Range("J118").Value = 2800
$ws.Range("L118").Value = 8400
$ws.Range("N118").Value = -10886

$ws.Range("H122").Value = 836.1818
$ws.Range("I122").Value = 417.14285
$ws.Range("J122").Value = 1569.5
$ws.Range("K122").Value = 3754.28565
$ws.Range("L122").Value = 14125.5
$ws.Range("M122").Value = -1304.28565
$ws.Range("N122").Value = -19025.5

$ws.Range("H131").Value = 1395.551
$ws.Range("J131").Value = 1094.3
$ws.Range("L131").Value = 3282.9
$ws.Range("N131").Value = -13362.9

$ws.Range("H135").Value = 276.02777
$ws.Range("I135").Value = 276.02777
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2484.24993
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = 50.75007000000005
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 149083.5
$ws.Range("J18").Value = 52603.75
$ws.Range("L18").Value = 52603.75
$ws.Range("N18").Value = -53189.75

$ws.Range("H43").Value = 3350.75
$ws.Range("I43").Value = 2746.75
$ws.Range("J43").Value = 3954.75
$ws.Range("K43").Value = 2746.75
$ws.Range("L43").Value = 3954.75
$ws.Range("M43").Value = -2595.75
$ws.Range("N43").Value = -4256.75

$ws.Range("H46").Value = 6569.2
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 8923
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 8923
$ws.Range("M46").Value = -4844
$ws.Range("N46").Value = -9235

$ws.Range("H57").Value = 11843.5
$ws.Range("I57").Value = 11200
$ws.Range("J57").Value = 15061
$ws.Range("K57").Value = 11200
$ws.Range("L57").Value = 15061
$ws.Range("M57").Value = -10380
$ws.Range("N57").Value = -16701

$ws.Range("H80").Value = 3118.182
$ws.Range("I80").Value = 2666.6667
$ws.Range("J80").Value = 3660
$ws.Range("K80").Value = 2666.6667
$ws.Range("L80").Value = 3660
$ws.Range("M80").Value = -1668.6667
$ws.Range("N80").Value = -5656

$ws.Range("H83").Value = 3118.182
$ws.Range("I83").Value = 2666.6667
$ws.Range("J83").Value = 3660
$ws.Range("K83").Value = 13333.3335
$ws.Range("L83").Value = 18300
$ws.Range("M83").Value = -8341.333500000001
$ws.Range("N83").Value = -28284

$ws.Range("H102").Value = 21504.365
$ws.Range("I102").Value = 1467.5
$ws.Range("J102").Value = 44880.707
$ws.Range("K102").Value = 1467.5
$ws.Range("L102").Value = 44880.707
$ws.Range("M102").Value = 154.5
$ws.Range("N102").Value = -48124.707

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1730.8
$ws.Range("I7").Value = 1413.5
$ws.Range("J7").Value = 1942.3334
$ws.Range("K7").Value = 1413.5
$ws.Range("L7").Value = 1942.3334
$ws.Range("M7").Value = -1301.5
$ws.Range("N7").Value = -2166.3334

$ws.Range("H18").Value = 56668.668
$ws.Range("I18").Value = 10000
$ws.Range("K18").Value = 10000
$ws.Range("M18").Value = -9828

$ws.Range("H20").Value = 50337.332
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 50337.332
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 50337.332
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -50789.332

$ws.Range("H126").Value = 1730.8
$ws.Range("I126").Value = 1413.5
$ws.Range("J126").Value = 1942.3334
$ws.Range("K126").Value = 4240.5
$ws.Range("L126").Value = 5827.0002
$ws.Range("M126").Value = -1770.5
$ws.Range("N126").Value = -10767.0002

$ws.Range("H136").Value = 4381.5
$ws.Range("I136").Value = 4282.5713
$ws.Range("J136").Value = 4520
$ws.Range("K136").Value = 12847.7139
$ws.Range("L136").Value = 13560
$ws.Range("M136").Value = -10297.7139
$ws.Range("N136").Value = -18660
